# LocalAirData.xlsx test-data refresh:
#  - ToLocation for the DataProviderWithExcel_002 row changes from "bost" to "miami"
#  - the sheet's remembered selection/scroll moves from D11 (scrolled to show column B)
#    to D9 (no scroll offset)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "miami"

$ws.Range("D9").Select() | Out-Null
